$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 320, pushing existing rows 320-374 down to 322-376
$ws.Rows("320:321").Insert()

# Populate the two newly inserted rows (320 and 321) with the new data entries
$ws.Range("A320").Value = 3
$ws.Range("B320").Value = "Femacal de La Calera"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = 44776
$ws.Range("E320").Value = 5
$ws.Range("F320").Value = 100112013
$ws.Range("G320").Value = "Alcachofa"
$ws.Range("H320").Value = "Argentina(o)"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 130
$ws.Range("K320").Value = 14500
$ws.Range("L320").Value = 15000
$ws.Range("M320").Value = 14769
$ws.Range("N320").Value = "$/caja 50 unidades"
$ws.Range("O320").Value = "Provincia de Limarí"
$ws.Range("P320").Value = 295
$ws.Range("Q320").Value = 50
$ws.Range("R320").Value = "Hortaliza"

$ws.Range("A321").Value = 3
$ws.Range("B321").Value = "Femacal de La Calera"
$ws.Range("C321").Value = "Coquimbo"
$ws.Range("D321").Value = 44776
$ws.Range("E321").Value = 5
$ws.Range("F321").Value = 100112013
$ws.Range("G321").Value = "Alcachofa"
$ws.Range("H321").Value = "Española"
$ws.Range("I321").Value = "Extra"
$ws.Range("J321").Value = 115
$ws.Range("K321").Value = 15000
$ws.Range("L321").Value = 16000
$ws.Range("M321").Value = 15522
$ws.Range("N321").Value = "$/caja 30 unidades"
$ws.Range("O321").Value = "Provincia de Limarí"
$ws.Range("P321").Value = 517
$ws.Range("Q321").Value = 30
$ws.Range("R321").Value = "Hortaliza"
